$d = $word.ActiveDocument

# --- Edit 1: underline "loguearme" in the first user story -----------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("loguearme", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $find1.Font.Underline = 1
}

# --- Edit 2: add user story 13 ("registrarme") ------------------------------
# Locate the existing "asignar una puntuación a las atracciones turísticas"
# paragraph (the last numbered user-story item already in the document) so
# the new story paragraph can be appended right after it, inheriting the
# same numbered-list paragraph formatting.
$anchorRange = $d.Content
$anchorRange.Find.Execute("asignar una puntuación a las atracciones turísticas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchorRange.Paragraphs.Item(1)

$anchorPara.Range.InsertParagraphAfter()
$newIndex = $anchorPara.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.Text = "Yo como usuario deseo poder registrarme en la página, por medio de un usuario y contraseña."

# Underline "usuario" and "registrarme" the same way the rest of the
# document marks the actor / action words of each user story.
$storyRange = $newPara.Range

$underline1 = $storyRange.Duplicate
$underline1.Find.Execute("usuario ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$underline1.Font.Underline = 1

$underline2 = $storyRange.Duplicate
$underline2.Find.Execute("registrarme", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$underline2.Font.Underline = 1

# Add the trailing blank paragraph that separates this story from the next
# block, matching the other blank separator already in the document: no
# list numbering, indented to match the list body (0.5in / 720 twips).
$newPara.Range.InsertParagraphAfter()
$blankIndex = $newIndex + 1
$blankPara = $d.Paragraphs.Item($blankIndex)

$blankXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>
'@
$blankPara.Range.InsertXML($blankXml)
